$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.413.91"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.863.75"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7043"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9989"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07770"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3086"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07837"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.172"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "93.50"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "1.859.86"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.643"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008377"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "29.381.17"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "2.098.70"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.591"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9986"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1523"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.939"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.45"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.257"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.224"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.202"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05160"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7941"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.932"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.48%  "
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.691"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").Value = "1.334.94"
$ws.Range("E38").Value = "  +8.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01879"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.734"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9593"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.082"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.72"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9987"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.809"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("D46").Value = "1.999.63"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.26"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5194"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("E50").Value = "  -5.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.031"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.34%  "
